$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '57651'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2625'
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = '64109'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2520'
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = '45624'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3083'
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '61342'
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '9677'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '49710892'
$ws.Range("C8").Value = 'MMMMMMM'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4746'
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '11780'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '53060417'
$ws.Range("C9").Value = '㊥老纳信耶稣'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4640'
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '15035'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4491'
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '16914'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4418'
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '19619'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '4319'
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '39369'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3542'
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '58653'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2605'
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = '65937'
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '11586'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '4650'
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '12775'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '31495601'
$ws.Range("C18").Value = '陈晓军'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '4591'
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = '13272'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '54698813'
$ws.Range("C19").Value = '閃亮唐老鴨'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '4567'
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = '20919'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '54085771'
$ws.Range("C20").Value = '㊥Matthieu'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4276'
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = '21558'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '46289694'
$ws.Range("C21").Value = '㊥Vincent'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4255'
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = '21858'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '55769051'
$ws.Range("C22").Value = '㊥叮叮当.'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4245'
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = '21882'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '8057001'
$ws.Range("C23").Value = '㊥兵者诡道也'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4244'
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = '26545'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '4101'
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = '27001'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '3649043'
$ws.Range("C25").Value = 'Dj6106'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '4089'
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = '29768'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4006'
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = '30265'
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '58839983'
$ws.Range("C27").Value = '每逢佳节胖六斤'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '3996'
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = '31929'
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = '58408326'
$ws.Range("C28").Value = '"Killer Bee"'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '3983'
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = '40271'
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = '1304123'
$ws.Range("C29").Value = 'Cccccccccccc'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '3477'
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = '40644'
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '55860890'
$ws.Range("C30").Value = '㊥Ethan'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '3446'
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = '1355'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '5340'
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = '6629'
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = '35114520'
$ws.Range("C32").Value = '13lur¹³'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '4907'
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = '8092'
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = '7852598'
$ws.Range("C33").Value = 'seiji'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '4827'
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = '9591'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '4752'
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = '10309'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '4714'
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = '12681'
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = '55317038'
$ws.Range("C36").Value = 'necman12345'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '4595'
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = '14871'
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = '26588375'
$ws.Range("C37").Value = '何苦僧ᶻᵍˣ'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '4497'
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = '15971'
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = '26280580'
$ws.Range("C38").Value = '꧁SSS.TIGRESS꧂ᶻᵍˣ'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '4454'
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = '18753'
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = '56133764'
$ws.Range("C39").Value = 'ustcarter'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4350'
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = '20515'
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = '6809364'
$ws.Range("C40").Value = '"Scorp IP"'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4288'
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = '22230'
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = '38995116'
$ws.Range("C41").Value = '"Ramesh Pavai Nam"'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4233'
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = '26140'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4113'
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = '28015'
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = '56573048'
$ws.Range("C43").Value = 'Xiaotian'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4059'
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = '28269'
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = '56379103'
$ws.Range("C44").Value = 'Globalking'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4051'
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = '32887'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3954'
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = '36360'
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = '58203298'
$ws.Range("C46").Value = '权旨qua'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3751'
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = '38326'
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = '52997727'
$ws.Range("C47").Value = 'larios'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '3619'
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = '41070'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3413'
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = '41931'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '3350'
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = '52611'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '2766'
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = '54637'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '2702'
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = '58443'
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = '60309'
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = '66059'
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = '59348'
$ws.Range("E55").NumberFormat = "@"
$ws.Range("E55").Value = '2591'
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = '56941'
$ws.Range("E56").NumberFormat = "@"
$ws.Range("E56").Value = '2641'
$ws.Range("A59").NumberFormat = "@"
$ws.Range("A59").Value = '31553'
$ws.Range("E59").NumberFormat = "@"
$ws.Range("E59").Value = '3988'
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = '44861'
$ws.Range("E60").NumberFormat = "@"
$ws.Range("E60").Value = '3132'
$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = '47845'
$ws.Range("E61").NumberFormat = "@"
$ws.Range("E61").Value = '2958'
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = '58470'
$ws.Range("E62").NumberFormat = "@"
$ws.Range("E62").Value = '2608'
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = '65383'
$ws.Range("E63").NumberFormat = "@"
$ws.Range("E63").Value = '2505'
$ws.Range("A64").NumberFormat = "@"
$ws.Range("A64").Value = '95975'
$ws.Range("E64").NumberFormat = "@"
$ws.Range("E64").Value = '1508'
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = '98869'
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = '15436348'
$ws.Range("C65").Value = 'Lucas'
$ws.Range("E65").NumberFormat = "@"
$ws.Range("E65").Value = '1497'
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = '111486'
$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = '49000199'
$ws.Range("C66").Value = 'SlipperyForester5672'
$ws.Range("E66").NumberFormat = "@"
$ws.Range("E66").Value = '1278'
$ws.Range("B67").NumberFormat = "@"
$ws.Range("B67").Value = '56700848'
$ws.Range("C67").Value = '工口漫画老师'
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = '38994054'
$ws.Range("C68").Value = 'chengnan'
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = '3391765'
$ws.Range("C69").Value = '马er'
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = '55810157'
$ws.Range("C70").Value = 'Beard'
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = '57556179'
$ws.Range("C71").Value = '特战新生代英雄'
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = '1222440'
$ws.Range("C72").Value = '"Sneaky Ninja Panda"'
$ws.Range("B73").NumberFormat = "@"
$ws.Range("B73").Value = '58340439'
$ws.Range("C73").Value = '70qilin'
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = '51743'
$ws.Range("E79").NumberFormat = "@"
$ws.Range("E79").Value = '2795'
$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = '117592'
$ws.Range("E82").NumberFormat = "@"
$ws.Range("E82").Value = '1187'
$ws.Range("A83").NumberFormat = "@"
$ws.Range("A83").Value = '134239'
$ws.Range("A84").NumberFormat = "@"
$ws.Range("A84").Value = '169658'
